$d = $word.ActiveDocument

$d.Content.Find.Execute("2026-01-28 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2026-01-29 Thursday", 2) | Out-Null
$d.Content.Find.Execute("951×2=", $true, $false, $false, $false, $false, $true, 1, $false, "512×3=", 2) | Out-Null
$d.Content.Find.Execute("817×5=", $true, $false, $false, $false, $false, $true, 1, $false, "221×2=", 2) | Out-Null
$d.Content.Find.Execute("146×2=", $true, $false, $false, $false, $false, $true, 1, $false, "206×4=", 2) | Out-Null
$d.Content.Find.Execute("938×2=", $true, $false, $false, $false, $false, $true, 1, $false, "543×4=", 2) | Out-Null
$d.Content.Find.Execute("190×2=", $true, $false, $false, $false, $false, $true, 1, $false, "775×8=", 2) | Out-Null
$d.Content.Find.Execute("244×6=", $true, $false, $false, $false, $false, $true, 1, $false, "991×5=", 2) | Out-Null
$d.Content.Find.Execute("182×7=", $true, $false, $false, $false, $false, $true, 1, $false, "730×6=", 2) | Out-Null
$d.Content.Find.Execute("861×8=", $true, $false, $false, $false, $false, $true, 1, $false, "655×4=", 2) | Out-Null
$d.Content.Find.Execute("203×4=", $true, $false, $false, $false, $false, $true, 1, $false, "333×7=", 2) | Out-Null
$d.Content.Find.Execute("103×6=", $true, $false, $false, $false, $false, $true, 1, $false, "648×8=", 2) | Out-Null
$d.Content.Find.Execute("509×4=", $true, $false, $false, $false, $false, $true, 1, $false, "257×7=", 2) | Out-Null
$d.Content.Find.Execute("827×4=", $true, $false, $false, $false, $false, $true, 1, $false, "497×8=", 2) | Out-Null
$d.Content.Find.Execute("919×9=", $true, $false, $false, $false, $false, $true, 1, $false, "604×7=", 2) | Out-Null
$d.Content.Find.Execute("273×2=", $true, $false, $false, $false, $false, $true, 1, $false, "739×3=", 2) | Out-Null
$d.Content.Find.Execute("235×7=", $true, $false, $false, $false, $false, $true, 1, $false, "436×4=", 2) | Out-Null
$d.Content.Find.Execute("693×2=", $true, $false, $false, $false, $false, $true, 1, $false, "451×8=", 2) | Out-Null
$d.Content.Find.Execute("604×8=", $true, $false, $false, $false, $false, $true, 1, $false, "986×7=", 2) | Out-Null
$d.Content.Find.Execute("684×2=", $true, $false, $false, $false, $false, $true, 1, $false, "659×8=", 2) | Out-Null
$d.Content.Find.Execute("376×4=", $true, $false, $false, $false, $false, $true, 1, $false, "530×6=", 2) | Out-Null
$d.Content.Find.Execute("618×9=", $true, $false, $false, $false, $false, $true, 1, $false, "738×3=", 2) | Out-Null
$d.Content.Find.Execute("753×7=", $true, $false, $false, $false, $false, $true, 1, $false, "346×8=", 2) | Out-Null
$d.Content.Find.Execute("620×7=", $true, $false, $false, $false, $false, $true, 1, $false, "988×6=", 2) | Out-Null
$d.Content.Find.Execute("300×3=", $true, $false, $false, $false, $false, $true, 1, $false, "721×3=", 2) | Out-Null
$d.Content.Find.Execute("960×7=", $true, $false, $false, $false, $false, $true, 1, $false, "608×6=", 2) | Out-Null
$d.Content.Find.Execute("626×6=", $true, $false, $false, $false, $false, $true, 1, $false, "254×6=", 2) | Out-Null
